# Updates the cryptos price/volume table (and re-orders two coin pairs)
# to reflect the latest scrape, per the GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.785.23'
$ws.Range("E2").Value = '  -1.14%  '
$ws.Range("D3").Value = '2.239.69'
$ws.Range("E3").Value = '  -1.85%  '
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '113.54'
$ws.Range("E5").Value = '  -0.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '269.39'
$ws.Range("E6").Value = '  +1.48%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.627'
$ws.Range("E7").Value = '  +1.61%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.606'
$ws.Range("E9").Value = '  -0.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '46.16'
$ws.Range("E10").Value = '  -2.41%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0930'
$ws.Range("E11").Value = '  -0.51%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.12'
$ws.Range("E12").Value = '  -1.62%  '
$ws.Range("E13").Value = '  -2.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.31'
$ws.Range("E14").Value = '  -1.05%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.873'
$ws.Range("E15").Value = '  +0.72%  '
$ws.Range("D16").Value = '2.575.78'
$ws.Range("E16").Value = '  -1.69%  '
$ws.Range("D17").Value = '2.241.18'
$ws.Range("E17").Value = '  -1.68%  '
$ws.Range("D18").Value = '42.977.11'
$ws.Range("E18").Value = '  -0.64%  '
$ws.Range("E19").Value = '  -1.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.75'
$ws.Range("E20").Value = '  -0.61%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.93'
$ws.Range("E21").Value = '  -0.06%  '
$ws.Range("E22").Value = '  -5.73%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '230.93'
$ws.Range("E23").Value = '  -1.27%  '
$ws.Range("E24").Value = '  +2.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.31'
$ws.Range("E25").Value = '  -3.22%  '
$ws.Range("E26").Value = '  +7.20%  '
$ws.Range("E27").Value = '  -0.35%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '40.28'
$ws.Range("E28").Value = '  -2.05%  '
$ws.Range("E29").Value = '  -1.00%  '
$ws.Range("E30").Value = '  -1.78%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '173.62'
$ws.Range("E31").Value = '  -0.11%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.10'
$ws.Range("E32").Value = '  -2.58%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0901'
$ws.Range("E33").Value = '  -0.62%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.57'
$ws.Range("E34").Value = '  -2.63%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.27'
$ws.Range("E35").Value = '  +8.88%  '
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.71'
$ws.Range("E37").Value = '  +0.97%  '
$ws.Range("E38").Value = '  +1.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.106'
$ws.Range("E39").Value = '  +3.41%  '
$ws.Range("E40").Value = '  -2.42%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '70.87'
$ws.Range("E41").Value = '  -7.46%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '13.15'
$ws.Range("E42").Value = '  -7.37%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.232'
$ws.Range("E43").Value = '  -2.73%  '
$ws.Range("E44").Value = '  +0.32%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.66'
$ws.Range("E45").Value = '  -8.65%  '
$ws.Range("E46").Value = '  -3.67%  '
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.45'
$ws.Range("E47").Value = '  -1.57%  '
$ws.Range("B48").Value = 'TrustWalletToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.25'
$ws.Range("E48").Value = '  -0.46%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0987'
$ws.Range("E49").Value = '  -0.93%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '100.24'
$ws.Range("E50").Value = '  -3.40%  '
$ws.Range("B51").Value = 'TheSandbox'
$ws.Range("C51").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.641'
$ws.Range("E51").Value = '  +7.17%  '
